$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 (Stocks / "Third Manual Column Reduce...") -- List-of-Cols (D) gains 4 new columns
$ws.Range("D9").Value = "['lco', 'lcox', 'lcoxdr', 'lo', 'loxdr', 'mib', 'mibn', 'ppeveb', 'pstkc', 'pstkl', 'pstkn', 'pstkr', 'tstkc', 'tstkn', 'che', 'cicurr', 'cidergl', 'cimii', 'ciother', 'cipen', 'cisecgl', 'citotal', 'dpc', 'dpvieb', 'dv', 'dvp', 'dvpa', 'ib', 'epspx', 'esopnr', 'esopt', 'ibadj', 'ibc', 'ibcom', 'ibmii', 'recch', 'recco', 'rectr', 'reuna', 'sale', 'spced', 'spceeps', 'cshtr_c', 'dvpsp_c', 'dvpsx_c', 'prcc_c', 'prch_c', 'prcl_c', 'adjex_c', 'acdo', 'aco', 'acodo', 'acox', 'aldo', 'aocidergl', 'aociother', 'aocipen', 'aodo', 'aox', 'ap', 'ceql', 'intc', 'ivaco', 'ivaeq', 'ivao', 'ivch', 'ivst', 'ivstch', 'pncad', 'pncaeps', 'prcad', 'prcaeps', 'xido', 'xidoc', 'ajex', 'ajp', 'cshfd', 'cshi', 'csho', 'cstk', 'cstkcv', 'cstke', 'dclo', 'dcom', 'dcvsr', 'dcvsub', 'dcvt', 'dd', 'dd1', 'dd2', 'dltis', 'dlto', 'dm', 'dn', 'ds', 'dudd', 'fatc', 'fatc', 'fatn', 'fiao', 'fopox', 'intano', 'mrc1', 'mrcta', 'niadj', 'nopio', 'oiadp', 'oibdp', 'oprepsx', 'pnrsho', 'prsho', 'ppent', 'pstkrv', 'txbco', 'txbcof', 'txdba', 'txdbca', 'txdbcl', 'txdc', 'txdi', 'txditc', 'txndba', 'txndb', 'txndbl', 'txo', 'txp', 'txpd', 'txr', 'acctstd', 'am', 'capxv', 'dc', 'diladj', 'do', 'donr', 'emp', 'esub', 'exre', 'lifr', 'mibt', 'prstkc', 'seqo', 'spi', 'cshtr_f', 'dvpsp_f', 'dvpsx_f', 'prcc_f', 'prch_f', 'prcl_f', 'adjex_f', 'ebitda', 'acominc', 'act', 'ao', 'aqc']"

# Row 28 (Stocks / "Manual Drop -- 2...") -- List-of-Cols (D) gains 'cshoc'
$ws.Range("D28").Value = "['iid', 'exchg', 'tpci', 'cik', 'cshtrd', 'ajexdi', 'trfd', 'cshoc']"

# ResultShape (E) updates reflecting the new, slightly smaller column counts
$ws.Range("E9").Value = "(1243, 86)"
$ws.Range("E10").Value = "(1243, 79)"
$ws.Range("E11").Value = "(348, 149)"
$ws.Range("E12").Value = "(348, 154)"
$ws.Range("E28").Value = "(673564, 5)"
$ws.Range("E29").Value = "(689, 7)"
$ws.Range("E30").Value = "(689, 6)"
$ws.Range("E35").Value = "(348, 154)"
$ws.Range("E36").Value = "(348, 159)"
$ws.Range("E37").Value = "(348, 165)"
$ws.Range("E38").Value = "(348, 166)"
$ws.Range("E39").Value = "(348, 169)"
